$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) cells for each coin row.
# NumberFormat is forced to Text ("@") before assigning values that look like
# plain numbers (e.g. "1.00", "168.52") so Excel keeps them as literal text
# instead of re-parsing/rounding them as numeric values; ClearFormats()
# afterwards removes the temporary formatting so the cell style is unchanged.

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "67.013.99"
$cell.ClearFormats()
$ws.Cells.Item(2, 5).Value = "  -0.04%  "

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.078.93"
$cell.ClearFormats()
$ws.Cells.Item(3, 5).Value = "  -1.12%  "

$ws.Cells.Item(4, 5).Value = "  +0.29%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "576.06"
$cell.ClearFormats()
$ws.Cells.Item(5, 5).Value = "  -0.13%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "168.52"
$cell.ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -3.04%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Cells.Item(7, 5).Value = "  +0.07%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.078.87"
$cell.ClearFormats()
$ws.Cells.Item(8, 5).Value = "  -0.94%  "

$ws.Cells.Item(9, 5).Value = "  -1.23%  "

$ws.Cells.Item(10, 5).Value = "  -0.79%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.150"
$cell.ClearFormats()
$ws.Cells.Item(11, 5).Value = "  -1.73%  "

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.471"
$cell.ClearFormats()
$ws.Cells.Item(12, 5).Value = "  -1.66%  "

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000241"
$cell.ClearFormats()
$ws.Cells.Item(13, 5).Value = "  -2.03%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "36.08"
$cell.ClearFormats()
$ws.Cells.Item(14, 5).Value = "  -2.87%  "

$ws.Cells.Item(15, 5).Value = "  -2.07%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.597.84"
$cell.ClearFormats()
$ws.Cells.Item(16, 5).Value = "  -0.77%  "

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = "67.041.75"
$cell.ClearFormats()
$ws.Cells.Item(17, 5).Value = "  +0.12%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.01"
$cell.ClearFormats()
$ws.Cells.Item(18, 5).Value = "  -1.40%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.66"
$cell.ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +2.25%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.084.80"
$cell.ClearFormats()
$ws.Cells.Item(20, 5).Value = "  -0.82%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "491.80"
$cell.ClearFormats()
$ws.Cells.Item(21, 5).Value = "  +3.01%  "

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.687"
$cell.ClearFormats()
$ws.Cells.Item(22, 5).Value = "  -3.64%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.69"
$cell.ClearFormats()
$ws.Cells.Item(23, 5).Value = "  -2.36%  "

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "82.88"
$cell.ClearFormats()
$ws.Cells.Item(24, 5).Value = "  -1.12%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.86"
$cell.ClearFormats()
$ws.Cells.Item(25, 5).Value = "  -4.31%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.21"
$cell.ClearFormats()
$ws.Cells.Item(26, 5).Value = "  -3.26%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.22"
$cell.ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +2.93%  "

$ws.Cells.Item(28, 5).Value = "  -0.01%  "

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.83"
$cell.ClearFormats()
$ws.Cells.Item(29, 5).Value = "  -1.34%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.29"
$cell.ClearFormats()
$ws.Cells.Item(30, 5).Value = "  -6.20%  "

$ws.Cells.Item(31, 5).Value = "  -1.52%  "

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "27.85"
$cell.ClearFormats()
$ws.Cells.Item(32, 5).Value = "  -2.94%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.111"
$cell.ClearFormats()
$ws.Cells.Item(33, 5).Value = "  -2.26%  "

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0905"
$cell.ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -4.28%  "

$ws.Cells.Item(35, 5).Value = "  -0.01%  "

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.68"
$cell.ClearFormats()
$ws.Cells.Item(36, 5).Value = "  -3.06%  "

$ws.Cells.Item(37, 5).Value = "  -2.57%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "46.87"
$cell.ClearFormats()
$ws.Cells.Item(38, 5).Value = "  -1.71%  "

$ws.Cells.Item(39, 5).Value = "  +0.79%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.99"
$cell.ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -4.88%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.302"
$cell.ClearFormats()
$ws.Cells.Item(41, 5).Value = "  -2.51%  "

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.31"
$cell.ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -3.46%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.771.02"
$cell.ClearFormats()
$ws.Cells.Item(43, 5).Value = "  -0.91%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "372.25"
$cell.ClearFormats()
$ws.Cells.Item(44, 5).Value = "  -1.88%  "

$ws.Cells.Item(45, 5).Value = "  -2.83%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "135.25"
$cell.ClearFormats()
$ws.Cells.Item(46, 5).Value = "  -0.66%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.45"
$cell.ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -4.84%  "

$ws.Cells.Item(48, 5).Value = "  +0.03%  "

$ws.Cells.Item(49, 5).Value = "  -1.49%  "

$ws.Cells.Item(50, 5).Value = "  -2.01%  "

$ws.Cells.Item(51, 5).Value = "  -1.59%  "
